$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Unhide all data rows first (writing new values into a hidden row
#    causes this engine to stamp an odd auto-fit height on it, so make
#    sure nothing is hidden before we touch any cell contents).
$ws.Range("A1:H23").EntireRow.Hidden = $false

# 2) New headers for the two additional columns.
$ws.Range("G1").Value = "Area (m^2)"
$ws.Range("H1").Value = "Area (Ha)"

# 3) New column data (G = Area in m^2, H = Area in Ha).
$areaM2 = @(5932.7340000000004,9214.1880000000001,9787.143,5911.09,6042.2740000000003,6122.893,11559.834000000001,9116.0849999999991,9358.5400000000009,6912.3890000000001,9240.3279999999995,7074.09,10156.235000000001,14011.459000000001,12328.22,11575.291999999999,13557.700999999999,7339.2049999999999,12526.308000000001,15052.763000000001,11011.726000000001,16063.66)
$areaHa = @(0.59299999999999997,0.92100000000000004,0.97899999999999998,0.59099999999999997,0.60399999999999998,0.61199999999999999,1.1559999999999999,0.91200000000000003,0.93600000000000005,0.69099999999999995,0.92400000000000004,0.70699999999999996,1.016,1.401,1.2330000000000001,1.1579999999999999,1.3560000000000001,0.73399999999999999,1.2529999999999999,1.5049999999999999,1.101,1.6060000000000001)

for ($i = 0; $i -lt 22; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 7).Value = $areaM2[$i]
    $ws.Cells.Item($r, 8).Value = $areaHa[$i]
}

# 4) Rebuild the table ("Tabla1") over the new A1:H23 range: this both
#    widens it to include the two new columns and drops the existing
#    "split = train" AutoFilter criteria (the rows are no longer
#    filtered, matching the new, unfiltered table).
$lo = $ws.ListObjects.Item(1)
$tableName = $lo.Name
$lo.Unlist()
$newLo = $ws.ListObjects.Add(1, $ws.Range("A1:H23"), 0, 1)
$newLo.Name = $tableName

# 5) Update the active selection to reflect where editing continued.
$ws.Range("I13").Select() | Out-Null
